$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "A: ['discounted_order'] == True"
$ws.Range("D2").Value = 1702.3519
$ws.Range("E2").Value = 97.0805
$ws.Range("F2").Value = 1889.9565

$ws.Range("C3").Value = "B: ['discounted_order'] == False"
$ws.Range("D3").Value = 1375.3318
$ws.Range("E3").Value = 84.5125
$ws.Range("F3").Value = 1790.7881
